# Apply the "Append: 2026-01-22 06:31 JST" scrape update to the
# "ランサーズ" (Lancers) listing sheet: the scraper re-ran, re-prioritised
# the list and only the top 3 rows survived the refresh (rows 5-13 from the
# previous run aged out of the window).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the old rows 5-13; this also shrinks the sheet dimension to A1:H4 ---
$ws.Rows("5:13").Delete()

# --- Row 2: MySQL / MariaDB database design & operations request ---
$ws.Range("A2").Value = "2026-01-22 06:31:53"
$ws.Range("B2").Value = "【急募】MySQL/MariaDBを活用したデータベース設計・運用の依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5476347"
$ws.Cells.Item(2, 7).Value = 33
$ws.Range("H2").Value = "◇MySQL"

# --- Row 3: long-term delivery-system Azure->AWS migration project ---
$ws.Range("A3").Value = "2026-01-22 06:31:53"
$ws.Range("B3").Value = "【長期案件】配送システム改修・AzureからAWS移行"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5476708"
$ws.Cells.Item(3, 7).Value = 40
$ws.Range("H3").ClearContents()

# --- Row 4: urgent private job posting ---
$ws.Range("A4").Value = "2026-01-22 06:31:53"
$ws.Range("B4").Value = "急募 限定公開 限定公開の仕事"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 60,000 円 / 募集期間 1 日、取引期間 0 日"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5476581"
$ws.Cells.Item(4, 7).Value = 13
$ws.Range("H4").ClearContents()

# --- Rebuild the hyperlinks: wipe the stale collection (rows 5-13's targets
# are now gone) and re-point F2:F4 at the 3 surviving/refreshed URLs ---
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5476347")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5476708")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5476581")
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"

# --- Column width adjustments (stored width = ColumnWidth + 5/6, so we
# subtract that offset to land on the exact target stored widths) ---
$ws.Columns("B").ColumnWidth = 37.166666666666664
$ws.Columns("D").ColumnWidth = 40.166666666666664
$ws.Columns("H").ColumnWidth = 11.166666666666666
